$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. This shifts the existing
# quantity/price/description/status columns one slot to the right
# (B->C, C->D, D->E, E->F).
$ws.Columns.Item(2).Insert()

# Rename the first header from "name" to "branch_name", and give the
# newly inserted column its own header. The inserted column already
# picked up the bold header style (s="1") from the row, so no extra
# formatting step is required.
$ws.Range("A1").Value = "branch_name"
$ws.Range("B1").Value = "product_name"

# Match the column widths implied by the target workbook (closest values
# reachable through the ColumnWidth character-width grid).
$ws.Columns.Item(1).ColumnWidth = 12.6
$ws.Columns.Item(2).ColumnWidth = 13.6

# Put the selection on B1, matching the saved view state.
$ws.Range("B1").Select()
